# Update Name of Algo
# Apply updated RandomForest imputation results to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value  = -13.28519999999999
$ws.Range("C4").Value  = -14.54600000000001
$ws.Range("C7").Value  = -11.6294
$ws.Range("C8").Value  = -12.36509999999999
$ws.Range("A11").Value = -21.85650000000002
$ws.Range("A12").Value = -21.45060000000002
$ws.Range("C12").Value = -11.92679999999999
$ws.Range("C14").Value = -11.49649999999999
$ws.Range("A15").Value = -21.26150000000001
$ws.Range("C22").Value = -10.98609999999999
